$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 282 ("08-10-2021") was previously missing its B/C values - fill them
# in now, matching the pattern used by every other row.
$ws.Range("B282").Value = 187
$ws.Range("C282").Value = 628

# Append the new daily rows (283-288), continuing the date series in column A
# and the repeating B/C/D/E pattern. The dates below (day <= 12) would
# otherwise be auto-converted to date serials by a plain .Value assignment,
# so they're written as a text formula first and then converted in place to
# a literal value via copy / paste-special (values only) - this keeps them
# stored as plain text, exactly like the rest of the column.
$newDates = "09-10-2021", "10-10-2021", "11-10-2021", "12-10-2021", "13-10-2021", "14-10-2021"

$row = 283
foreach ($date in $newDates) {
    $cellA = $ws.Range("A$row")
    $cellA.Formula = "=""" + $date + """"
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)

    # The final new row (288) keeps B/C empty, just like the old last row
    # (282) did before this update.
    if ($row -lt 288) {
        $ws.Range("B$row").Value = 187
        $ws.Range("C$row").Value = 628
    }
    $ws.Range("D$row").Value = 3940
    $ws.Range("E$row").Value = 30

    $row = $row + 1
}

$excel.CutCopyMode = $false
